$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 243; Time = "2023-12-12 16:19:07"; Cost = 0.0008 },
    @{ Row = 244; Time = "2023-12-12 16:19:29"; Cost = 0.001 },
    @{ Row = 245; Time = "2023-12-12 16:19:47"; Cost = 0.0012 },
    @{ Row = 246; Time = "2023-12-12 16:20:05"; Cost = 0.0004 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Time
    $ws.Cells.Item($r.Row, 2).Value = $r.Cost
}
